$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.422.29"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "1.848.35"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9986"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6324"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07565"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2955"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").Value = "1.853.11"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6862"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001007"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("D17").Value = "2.104.15"
$ws.Range("E17").Value = "  -1.29%  "
$ws.Range("D19").Value = "29.448.24"
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "228.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.539"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.13%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1396"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.375"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.60%  "
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05709"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.123"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.028"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.845"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.156"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7149"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.61%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "1.250.35"
$ws.Range("E38").Value = "  -0.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01809"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.777"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9098"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.176"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.091"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000117"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4025"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.122"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.685"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.11%  "
$ws.Range("E51").Value = "  -0.29%  "
